# Generate Report for Handoff
# Updates the "c0ec1678-eef8-48ac-b841-2c87645e88a6.md" row's handoff/generate
# timestamps across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) for row 6 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-09-02 20:48:34"

# --- zh-cn sheet: "Latest Handoff Datetime" (column H) for row 6 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-09-02 20:48:29"

# --- de-de sheet: "Latest Handoff Datetime" (column H) for row 6 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-09-02 20:48:34"
